# xls export geometry fix, added server power/health status
#
# Renames a handful of column headers to shorter labels (to match the new,
# narrower column geometry) and shrinks the corresponding columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label text ---------------------------------------------------
$ws.Range("F1").Value = "Memory tot.size"   # was: System memory size
$ws.Range("H1").Value = "Memory P/Ns"       # was: Memory module part number
$ws.Range("M1").Value = "HDD slot pop."     # was: HDD slot population
$ws.Range("N1").Value = "PSU P/Ns"          # was: PSU part number

# --- Column geometry -------------------------------------------------------
# Target widths (Excel "characters" units) derived from the stored
# OOXML <col/> width attribute (width = Truncate((N*MDW+5)/MDW*256)/256,
# MDW = 7 for the workbook's default Calibri 11 font):
#   col F  -> 15.7109375  => N = 15
#   col H  -> 11.7109375  => N = 11
#   col M  -> 13.7109375  => N = 13
#   col N  ->  8.7109375  => N = 8
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(13).ColumnWidth = 13
$ws.Columns.Item(14).ColumnWidth = 8
